$d = $word.ActiveDocument

# Fix a typo/missing-space in the "turn function" explanation paragraph:
# "transformedby" -> "transformed by"
$d.Content.Find.Execute("transformedby", $true, $false, $false, $false, $false,
                         $true, 1, $false, "transformed by", 2)
